$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 947.76666
$ws.Cells.Item(98, 9).Value = 1030.4348
$ws.Cells.Item(98, 11).Value = 1030.4348
$ws.Cells.Item(98, 13).Value = 467.5652
$ws.Cells.Item(116, 8).Value = 8100.3125
$ws.Cells.Item(116, 9).Value = 11640.5
$ws.Cells.Item(116, 11).Value = 11640.5
$ws.Cells.Item(116, 13).Value = -8198.5
$ws.Cells.Item(122, 8).Value = 947.76666
$ws.Cells.Item(122, 9).Value = 1030.4348
$ws.Cells.Item(122, 11).Value = 3091.3044
$ws.Cells.Item(122, 13).Value = -641.3044
$ws.Cells.Item(137, 8).Value = 1366.9822
$ws.Cells.Item(137, 9).Value = 941.86206
$ws.Cells.Item(137, 10).Value = 1823.5927
$ws.Cells.Item(137, 11).Value = 2825.58618
$ws.Cells.Item(137, 12).Value = 5470.7781
$ws.Cells.Item(137, 13).Value = -275.5861800000002
$ws.Cells.Item(137, 14).Value = -10570.7781
$ws.Cells.Item(138, 8).Value = 3043.2715
$ws.Cells.Item(138, 9).Value = 1342.3226
$ws.Cells.Item(138, 10).Value = 4395.3076
$ws.Cells.Item(138, 11).Value = 4026.9678
$ws.Cells.Item(138, 12).Value = 13185.9228
$ws.Cells.Item(138, 13).Value = 1113.0322
$ws.Cells.Item(138, 14).Value = -23465.9228
$ws.Cells.Item(139, 8).Value = 35000
$ws.Cells.Item(139, 10).Value = 35000
$ws.Cells.Item(139, 12).Value = 35000
$ws.Cells.Item(139, 14).Value = -45280

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1440.1904
$ws.Cells.Item(2, 9).Value = 1625.4706
$ws.Cells.Item(2, 10).Value = 652.75
$ws.Cells.Item(2, 11).Value = 1625.4706
$ws.Cells.Item(2, 12).Value = 652.75
$ws.Cells.Item(2, 13).Value = -1512.4706
$ws.Cells.Item(2, 14).Value = -878.75
$ws.Cells.Item(32, 8).Value = 4384.9536
$ws.Cells.Item(32, 9).Value = 4029.6365
$ws.Cells.Item(32, 11).Value = 4029.6365
$ws.Cells.Item(32, 13).Value = -3742.6365
$ws.Cells.Item(45, 8).Value = 7263.0415
$ws.Cells.Item(45, 9).Value = 7745.409
$ws.Cells.Item(45, 11).Value = 7745.409
$ws.Cells.Item(45, 13).Value = -7368.409
$ws.Cells.Item(61, 8).Value = 3239.9644
$ws.Cells.Item(61, 9).Value = 4360.757
$ws.Cells.Item(61, 10).Value = 1057.3684
$ws.Cells.Item(61, 11).Value = 4360.757
$ws.Cells.Item(61, 12).Value = 1057.3684
$ws.Cells.Item(61, 13).Value = -4148.757
$ws.Cells.Item(61, 14).Value = -1481.3684
$ws.Cells.Item(116, 8).Value = 1440.1904
$ws.Cells.Item(116, 9).Value = 1625.4706
$ws.Cells.Item(116, 10).Value = 652.75
$ws.Cells.Item(116, 11).Value = 1625.4706
$ws.Cells.Item(116, 12).Value = 652.75
$ws.Cells.Item(116, 13).Value = 668.5293999999999
$ws.Cells.Item(116, 14).Value = -5240.75
$ws.Cells.Item(124, 8).Value = 29166.666
$ws.Cells.Item(124, 10).Value = 29166.666
$ws.Cells.Item(124, 12).Value = 29166.666
$ws.Cells.Item(124, 14).Value = -38986.666
$ws.Cells.Item(125, 8).Value = 49939.2
$ws.Cells.Item(125, 10).Value = 49939.2
$ws.Cells.Item(125, 12).Value = 49939.2
$ws.Cells.Item(125, 14).Value = -59779.2
$ws.Cells.Item(136, 8).Value = 3239.9644
$ws.Cells.Item(136, 9).Value = 4360.757
$ws.Cells.Item(136, 10).Value = 1057.3684
$ws.Cells.Item(136, 11).Value = 13082.271
$ws.Cells.Item(136, 12).Value = 3172.1052
$ws.Cells.Item(136, 13).Value = -10532.271
$ws.Cells.Item(136, 14).Value = -8272.1052

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1440.1904
$ws.Cells.Item(3, 9).Value = 1625.4706
$ws.Cells.Item(3, 10).Value = 652.75
$ws.Cells.Item(3, 11).Value = 1625.4706
$ws.Cells.Item(3, 12).Value = 652.75
$ws.Cells.Item(3, 13).Value = -1511.4706
$ws.Cells.Item(3, 14).Value = -880.75
$ws.Cells.Item(99, 8).Value = 125002660
$ws.Cells.Item(99, 9).Value = 166669250
$ws.Cells.Item(99, 10).Value = 2900
$ws.Cells.Item(99, 11).Value = 166669250
$ws.Cells.Item(99, 12).Value = 2900
$ws.Cells.Item(99, 13).Value = -166667752
$ws.Cells.Item(99, 14).Value = -5896
$ws.Cells.Item(107, 8).Value = 895.9524
$ws.Cells.Item(107, 9).Value = 767.1539
$ws.Cells.Item(107, 10).Value = 1105.25
$ws.Cells.Item(107, 11).Value = 767.1539
$ws.Cells.Item(107, 12).Value = 1105.25
$ws.Cells.Item(107, 13).Value = 1152.8461
$ws.Cells.Item(107, 14).Value = -4945.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 70003
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(99, 8).Value = 12521601
$ws.Cells.Item(99, 9).Value = 23253
$ws.Cells.Item(99, 10).Value = 20853834
$ws.Cells.Item(99, 11).Value = 23253
$ws.Cells.Item(99, 12).Value = 20853834
$ws.Cells.Item(99, 13).Value = -21755
$ws.Cells.Item(99, 14).Value = -20856830
$ws.Cells.Item(109, 8).Value = 49150
$ws.Cells.Item(109, 10).Value = 49150
$ws.Cells.Item(109, 12).Value = 49150
$ws.Cells.Item(109, 14).Value = -51230
$ws.Cells.Item(126, 8).Value = 12521601
$ws.Cells.Item(126, 9).Value = 23253
$ws.Cells.Item(126, 10).Value = 20853834
$ws.Cells.Item(126, 11).Value = 69759
$ws.Cells.Item(126, 12).Value = 62561502
$ws.Cells.Item(126, 13).Value = -67289
$ws.Cells.Item(126, 14).Value = -62566442
$ws.Cells.Item(132, 8).Value = 4611.2856
$ws.Cells.Item(132, 9).Value = 2935.4285
$ws.Cells.Item(132, 10).Value = 6287.143
$ws.Cells.Item(132, 11).Value = 8806.2855
$ws.Cells.Item(132, 12).Value = 18861.429
$ws.Cells.Item(132, 13).Value = -6276.2855
$ws.Cells.Item(132, 14).Value = -23921.429

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 139482.48
$ws.Cells.Item(5, 10).Value = 167612.72
$ws.Cells.Item(5, 12).Value = 502838.16
$ws.Cells.Item(5, 14).Value = -503062.16
$ws.Cells.Item(68, 8).Value = 2614.3494
$ws.Cells.Item(68, 9).Value = 3562.2432
$ws.Cells.Item(68, 10).Value = 1851.9131
$ws.Cells.Item(68, 11).Value = 10686.7296
$ws.Cells.Item(68, 12).Value = 5555.7393
$ws.Cells.Item(68, 13).Value = -9875.729599999999
$ws.Cells.Item(68, 14).Value = -7177.7393
$ws.Cells.Item(71, 8).Value = 2614.3494
$ws.Cells.Item(71, 9).Value = 3562.2432
$ws.Cells.Item(71, 10).Value = 1851.9131
$ws.Cells.Item(71, 11).Value = 32060.1888
$ws.Cells.Item(71, 12).Value = 16667.2179
$ws.Cells.Item(71, 13).Value = -28004.1888
$ws.Cells.Item(71, 14).Value = -24779.2179
$ws.Cells.Item(132, 8).Value = 1236517.4
$ws.Cells.Item(132, 10).Value = 1524256.6
$ws.Cells.Item(132, 12).Value = 13718309.4
$ws.Cells.Item(132, 14).Value = -13723369.4
$ws.Cells.Item(135, 8).Value = 139482.48
$ws.Cells.Item(135, 10).Value = 167612.72
$ws.Cells.Item(135, 12).Value = 1508514.48
$ws.Cells.Item(135, 14).Value = -1513584.48
$ws.Cells.Item(137, 8).Value = 30318434
$ws.Cells.Item(137, 9).Value = 1514.625
$ws.Cells.Item(137, 10).Value = 40019850
$ws.Cells.Item(137, 11).Value = 4543.875
$ws.Cells.Item(137, 12).Value = 120059550
$ws.Cells.Item(137, 13).Value = 556.125
$ws.Cells.Item(137, 14).Value = -120069750

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6180.816
$ws.Cells.Item(70, 9).Value = 6297.5356
$ws.Cells.Item(70, 10).Value = 5854
$ws.Cells.Item(70, 11).Value = 6297.5356
$ws.Cells.Item(70, 12).Value = 5854
$ws.Cells.Item(70, 13).Value = -6027.5356
$ws.Cells.Item(70, 14).Value = -6394
$ws.Cells.Item(73, 8).Value = 6180.816
$ws.Cells.Item(73, 9).Value = 6297.5356
$ws.Cells.Item(73, 10).Value = 5854
$ws.Cells.Item(73, 11).Value = 6297.5356
$ws.Cells.Item(73, 12).Value = 5854
$ws.Cells.Item(73, 13).Value = -5361.5356
$ws.Cells.Item(73, 14).Value = -7726
$ws.Cells.Item(97, 8).Value = 820.6429000000001
$ws.Cells.Item(97, 9).Value = 861.25
$ws.Cells.Item(97, 10).Value = 766.5
$ws.Cells.Item(97, 11).Value = 861.25
$ws.Cells.Item(97, 12).Value = 766.5
$ws.Cells.Item(97, 13).Value = -365.25
$ws.Cells.Item(97, 14).Value = -1758.5
$ws.Cells.Item(140, 8).Value = 41781.58
$ws.Cells.Item(140, 10).Value = 41781.58
$ws.Cells.Item(140, 12).Value = 41781.58
$ws.Cells.Item(140, 14).Value = -52141.58

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 83336190
$ws.Cells.Item(40, 9).Value = 111113690
$ws.Cells.Item(40, 10).Value = 3711
$ws.Cells.Item(40, 11).Value = 111113690
$ws.Cells.Item(40, 12).Value = 3711
$ws.Cells.Item(40, 13).Value = -111113554
$ws.Cells.Item(40, 14).Value = -3983
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1741.4878
$ws.Cells.Item(132, 9).Value = 1598.7931
$ws.Cells.Item(132, 10).Value = 2086.3333
$ws.Cells.Item(132, 11).Value = 4796.379300000001
$ws.Cells.Item(132, 12).Value = 6258.999899999999
$ws.Cells.Item(132, 13).Value = -2266.379300000001
$ws.Cells.Item(132, 14).Value = -11318.9999
